$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new rows 32 and 33 with the Angular routing guard commands
$ws.Range("A32").Value = "Angular commands"
$ws.Range("B33").Value = "ng g guard auth --skipTests"
$ws.Range("A33").Value = "Generate guard"

# Update the view: select A33 and scroll so row 16 is the top-left cell
[void]$ws.Range("A33").Select()
$excel.ActiveWindow.ScrollRow = 16
